$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:AY that are always blank in the data rows being swapped (row 4/5/8/9/10/11).
# Copying a rectangular range stamps empty placeholder cells across the whole
# A:AY span; clearing these afterwards removes those placeholders again so the
# serialized row looks exactly like the source rows (sparse, no empty <c> tags).
$blankCols = @("C","J","K","L","M","N","O","X","Z","AB","AC","AF","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AU","AV")

$scratchRow = 100

function Clear-BlankCols($r) {
    foreach ($col in $blankCols) {
        $ws.Range("$col$r").ClearContents()
    }
}

function Swap-Rows($r1, $r2) {
    $src1 = "A$($r1):AY$($r1)"
    $src2 = "A$($r2):AY$($r2)"
    $scratch = "A$($scratchRow):AY$($scratchRow)"

    $ws.Range($src1).Copy($ws.Range($scratch))
    $ws.Range($src2).Copy($ws.Range($src1))
    $ws.Range($scratch).Copy($ws.Range($src2))
    $ws.Range($scratch).Clear()

    Clear-BlankCols $r1
    Clear-BlankCols $r2
}

Swap-Rows 4 5
Swap-Rows 8 9
Swap-Rows 10 11
